$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.485.32"
$ws.Range("D3").Value = "3.079.01"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.08"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.83"
$ws.Range("D8").Value = "3.073.36"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.57"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.35"
$ws.Range("D16").Value = "3.589.29"
$ws.Range("D17").Value = "66.558.67"
$ws.Range("D19").Value = "3.078.42"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "465.95"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.82"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.12"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "46.38"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "382.59"
$ws.Range("D47").Value = "2.752.23"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.02"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.57"

$ws.Range("E2").Value = "  +2.27%  "
$ws.Range("E3").Value = "  +4.35%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("E6").Value = "  +5.43%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +4.25%  "
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  +5.45%  "
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("E14").Value = "  +6.63%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("E16").Value = "  +4.38%  "
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("E20").Value = "  +14.44%  "
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("E22").Value = "  +4.76%  "
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("E25").Value = "  +5.03%  "
$ws.Range("E26").Value = "  +6.56%  "
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("E33").Value = "  +3.46%  "
$ws.Range("E34").Value = "  +4.17%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +2.46%  "
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("E38").Value = "  +7.58%  "
$ws.Range("E39").Value = "  +5.64%  "
$ws.Range("E40").Value = "  +6.33%  "
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E42").Value = "  +2.51%  "
$ws.Range("E43").Value = "  +3.04%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +3.65%  "
